$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1): extend the car ("汽車") sheet with the same
# trailing columns used by the other property sheets (property_category,
# category, date, legislator_name, legislator_id, source_file, index),
# plus renaming/adding the leading name/capacity/owner/... columns so the
# header matches the common record layout. Copy the existing header cell's
# formatting (bold font + border) onto the new header cells first so they
# look like the rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2): fill in the matching record values.
$ws.Range("A2").Value = 30
$ws.Range("B2").Value = "國瑞Wish"
$ws.Range("C2").Value = 1998
$ws.Range("D2").Value = "孫效智"
$ws.Range("E2").Value = "94年05月03H"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = "(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# "date" column: force text so the ISO-looking string isn't silently
# reinterpreted/reformatted as a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-24"

$ws.Range("K2").Value = "楊玉欣"
$ws.Range("L2").Value = 1757
$ws.Range("M2").Value = "tmp89971"
$ws.Range("N2").Value = 30
